$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 90 appended to the table (row count grows from 89 to 90 data rows).
$ws.Range("A90").Value = "Sagatova Ziyodaxon Taxirovna"
$ws.Range("B90").Value = "Maktabgacha talim tashkiloti tarbiyachisi"
$ws.Range("C90").Value = "AD1426517"

# D and G look numeric, but every other row in this sheet keeps them as plain
# text, so force text entry (leading apostrophe) and strip the resulting
# "number stored as text" quote-prefix style so the cell ends up unstyled,
# matching the rest of the sheet.
$ws.Range("D90").Value = "'268"
$ws.Range("D90").ClearFormats()

$ws.Range("E90").Value = "Toshkent viloyati"
$ws.Range("F90").Value = "Oʻrta Chirchiq tumani"

$ws.Range("G90").Value = "'998936160710"
$ws.Range("G90").ClearFormats()

$ws.Range("H90").Value = "23-11-2024"
